# Generate Report for Handback
# Update the "Latest/Correspond Handoff/Handback" datetime stamps that are
# (re)written when the handback report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$dtFormat = "yyyy-mm-dd HH:mm:ss"

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview.Range("G2").Value = "2016-08-31 11:11:21"
$wsOverview.Range("G2").NumberFormat = $dtFormat

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row
$wsZhCn.Range("H2").Value = "2016-08-31 11:11:16"
$wsZhCn.Range("H2").NumberFormat = $dtFormat
$wsZhCn.Range("K2").Value = "2016-08-31 11:11:49"
$wsZhCn.Range("K2").NumberFormat = $dtFormat

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row
$wsDeDe.Range("H2").Value = "2016-08-31 11:11:21"
$wsDeDe.Range("H2").NumberFormat = $dtFormat
$wsDeDe.Range("K2").Value = "2016-08-31 11:11:56"
$wsDeDe.Range("K2").NumberFormat = $dtFormat
